$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2017543859649123
$ws.Range("C2").Value = 0.5482456140350878
$ws.Range("J2").Value = 0.008771929824561403
$ws.Range("P2").Value = 0.1491228070175439
$ws.Range("S2").Value = 0.09210526315789473
$ws.Range("C3").Value = 0.01538461538461539
$ws.Range("J3").Value = 0.05384615384615385
$ws.Range("P3").Value = 0.7307692307692307
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.08163265306122448
$ws.Range("P4").Value = 0.6326530612244898
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("B6").Value = 0.05339805825242718
$ws.Range("D6").Value = 0.01941747572815534
$ws.Range("F6").Value = 0.05339805825242718
$ws.Range("J6").Value = 0.2864077669902912
$ws.Range("O6").Value = 0.02427184466019417
$ws.Range("Q6").Value = 0.1359223300970874
$ws.Range("R6").Value = 0.04854368932038835
$ws.Range("S6").Value = 0.3786407766990291
$ws.Range("B7").Value = 0.05027932960893855
$ws.Range("D7").Value = 0.03910614525139665
$ws.Range("E7").Value = 0.00558659217877095
$ws.Range("F7").Value = 0.07262569832402235
$ws.Range("J7").Value = 0.1620111731843575
$ws.Range("Q7").Value = 0.2011173184357542
$ws.Range("R7").Value = 0.07262569832402235
$ws.Range("S7").Value = 0.3966480446927375
$ws.Range("B8").Value = 0.07494145199063232
$ws.Range("D8").Value = 0.02576112412177986
$ws.Range("F8").Value = 0.04918032786885246
$ws.Range("J8").Value = 0.1030444964871194
$ws.Range("O8").Value = 0.01405152224824356
$ws.Range("Q8").Value = 0.1826697892271663
$ws.Range("R8").Value = 0.09836065573770492
$ws.Range("S8").Value = 0.4519906323185012
$ws.Range("B9").Value = 0.06072874493927125
$ws.Range("D9").Value = 0.01214574898785425
$ws.Range("F9").Value = 0.08097165991902834
$ws.Range("J9").Value = 0.1255060728744939
$ws.Range("O9").Value = 0.004048582995951417
$ws.Range("Q9").Value = 0.2024291497975708
$ws.Range("R9").Value = 0.08097165991902834
$ws.Range("S9").Value = 0.4331983805668016
$ws.Range("B10").Value = 0.08994276369582993
$ws.Range("D10").Value = 0.02207686017988553
$ws.Range("F10").Value = 0.06541291905151267
$ws.Range("J10").Value = 0.1095666394112837
$ws.Range("O10").Value = 0.01635322976287817
$ws.Range("Q10").Value = 0.2281275551921504
$ws.Range("R10").Value = 0.08340147179067865
$ws.Range("S10").Value = 0.3851185609157809
$ws.Range("G11").Value = 0.1618257261410788
$ws.Range("J11").Value = 0.07468879668049792
$ws.Range("K11").Value = 0.1867219917012448
$ws.Range("L11").Value = 0.5726141078838174
$ws.Range("S11").Value = 0.004149377593360996
$ws.Range("G12").Value = 0.7602739726027398
$ws.Range("J12").Value = 0.1849315068493151
$ws.Range("K12").Value = 0.00684931506849315
$ws.Range("L12").Value = 0.0136986301369863
$ws.Range("S12").Value = 0.03424657534246575
$ws.Range("G13").Value = 0.8222222222222222
$ws.Range("J13").Value = 0.1555555555555556
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.01913875598086124
$ws.Range("H15").Value = 0.1818181818181818
$ws.Range("I15").Value = 0.1004784688995215
$ws.Range("J15").Value = 0.354066985645933
$ws.Range("K15").Value = 0.05263157894736842
$ws.Range("M15").Value = 0.004784688995215311
$ws.Range("N15").Value = 0.004784688995215311
$ws.Range("O15").Value = 0.07177033492822966
$ws.Range("S15").Value = 0.2105263157894737
$ws.Range("F16").Value = 0.01935483870967742
$ws.Range("H16").Value = 0.1290322580645161
$ws.Range("I16").Value = 0.07741935483870968
$ws.Range("J16").Value = 0.4774193548387097
$ws.Range("K16").Value = 0.1354838709677419
$ws.Range("M16").Value = 0.01290322580645161
$ws.Range("N16").Value = 0.006451612903225806
$ws.Range("O16").Value = 0.03870967741935484
$ws.Range("S16").Value = 0.1032258064516129
$ws.Range("F17").Value = 0.01515151515151515
$ws.Range("H17").Value = 0.170995670995671
$ws.Range("I17").Value = 0.1233766233766234
$ws.Range("J17").Value = 0.4155844155844156
$ws.Range("K17").Value = 0.0735930735930736
$ws.Range("M17").Value = 0.01298701298701299
$ws.Range("N17").Value = 0.002164502164502165
$ws.Range("O17").Value = 0.0670995670995671
$ws.Range("S17").Value = 0.119047619047619
$ws.Range("F18").Value = 0.01621621621621622
$ws.Range("H18").Value = 0.145945945945946
$ws.Range("I18").Value = 0.1297297297297297
$ws.Range("J18").Value = 0.4432432432432433
$ws.Range("K18").Value = 0.06486486486486487
$ws.Range("M18").Value = 0.01621621621621622
$ws.Range("O18").Value = 0.05945945945945946
$ws.Range("S18").Value = 0.1243243243243243
$ws.Range("F19").Value = 0.01597444089456869
$ws.Range("H19").Value = 0.2100638977635783
$ws.Range("I19").Value = 0.1062300319488818
$ws.Range("J19").Value = 0.365814696485623
$ws.Range("K19").Value = 0.08945686900958466
$ws.Range("M19").Value = 0.02715654952076677
$ws.Range("O19").Value = 0.07108626198083066
$ws.Range("S19").Value = 0.1142172523961661
